$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintTracking")

$ws.Range("E3").Value = "Sarah"
$ws.Range("E4").Value = "Ani"
$ws.Range("E5").Value = "Nathan"

$ws.Range("C15").Select() | Out-Null
